# Scheduled runner update: refresh market-board derived leve profit figures
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 14287127
$ws.Range("I70").Value = 50000696
$ws.Range("K70").Value = 150002088
$ws.Range("M70").Value = -150001818

$ws.Range("H73").Value = 14287127
$ws.Range("I73").Value = 50000696
$ws.Range("K73").Value = 150002088
$ws.Range("M73").Value = -150001152

$ws.Range("H86").Value = 83371420
$ws.Range("I86").Value = 3773.5
$ws.Range("J86").Value = 125055250
$ws.Range("K86").Value = 3773.5
$ws.Range("L86").Value = 125055250
$ws.Range("M86").Value = -2650.5
$ws.Range("N86").Value = -125057496

$ws.Range("H89").Value = 83371420
$ws.Range("I89").Value = 3773.5
$ws.Range("J89").Value = 125055250
$ws.Range("K89").Value = 18867.5
$ws.Range("L89").Value = 625276250
$ws.Range("M89").Value = -13251.5
$ws.Range("N89").Value = -625287482

$ws.Range("H111").Value = 1153.5
$ws.Range("I111").Value = 1194.5
$ws.Range("K111").Value = 3583.5
$ws.Range("M111").Value = -516.5

$ws.Range("H112").Value = 2576.158
$ws.Range("I112").Value = 3511
$ws.Range("K112").Value = 10533
$ws.Range("M112").Value = -9425

$ws.Range("H116").Value = 12524551
$ws.Range("I116").Value = 12524551
$ws.Range("K116").Value = 12524551
$ws.Range("M116").Value = -12521109

$ws.Range("H138").Value = 260775.05
$ws.Range("I138").Value = 488496.28
$ws.Range("K138").Value = 1465488.84
$ws.Range("M138").Value = -1460348.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 69999.5
$ws.Range("J51").Value = 69999.5
$ws.Range("L51").Value = 69999.5
$ws.Range("N51").Value = -71511.5

$ws.Range("H63").Value = 2244.4666
$ws.Range("I63").Value = 2190.8572
$ws.Range("J63").Value = 2995
$ws.Range("K63").Value = 2190.8572
$ws.Range("L63").Value = 2995
$ws.Range("M63").Value = -1504.8572
$ws.Range("N63").Value = -4367

$ws.Range("H66").Value = 2244.4666
$ws.Range("I66").Value = 2190.8572
$ws.Range("J66").Value = 2995
$ws.Range("K66").Value = 10954.286
$ws.Range("L66").Value = 14975
$ws.Range("M66").Value = -7522.286
$ws.Range("N66").Value = -21839

$ws.Range("H92").Value = 183366670
$ws.Range("I92").Value = 90000
$ws.Range("J92").Value = 220022000
$ws.Range("K92").Value = 90000
$ws.Range("L92").Value = 220022000
$ws.Range("M92").Value = -87504
$ws.Range("N92").Value = -220026992

$ws.Range("H110").Value = 2149.1304
$ws.Range("I110").Value = 1715.4375
$ws.Range("K110").Value = 1715.4375
$ws.Range("M110").Value = 329.5625

$ws.Range("H122").Value = 837158.5600000001
$ws.Range("I122").Value = 3255.96
$ws.Range("K122").Value = 9767.880000000001
$ws.Range("M122").Value = -7317.880000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 35999.5
$ws.Range("J9").Value = 35999.5
$ws.Range("L9").Value = 35999.5
$ws.Range("N9").Value = -36335.5

$ws.Range("H64").Value = 14251.167
$ws.Range("I64").Value = 25666.666
$ws.Range("J64").Value = 2835.6667
$ws.Range("K64").Value = 25666.666
$ws.Range("L64").Value = 2835.6667
$ws.Range("M64").Value = -25441.666
$ws.Range("N64").Value = -3285.6667

$ws.Range("H67").Value = 14251.167
$ws.Range("I67").Value = 25666.666
$ws.Range("J67").Value = 2835.6667
$ws.Range("K67").Value = 25666.666
$ws.Range("L67").Value = 2835.6667
$ws.Range("M67").Value = -24886.666
$ws.Range("N67").Value = -4395.6667

$ws.Range("H134").Value = 8856.388999999999
$ws.Range("I134").Value = 9400.9375
$ws.Range("K134").Value = 28202.8125
$ws.Range("M134").Value = -25667.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1446.129
$ws.Range("I58").Value = 1200.88
$ws.Range("K58").Value = 1200.88
$ws.Range("M58").Value = -997.8800000000001

$ws.Range("H92").Value = 45583.332
$ws.Range("I92").Value = 40000
$ws.Range("J92").Value = 48375
$ws.Range("K92").Value = 40000
$ws.Range("L92").Value = 48375
$ws.Range("M92").Value = -37504
$ws.Range("N92").Value = -53367

$ws.Range("H136").Value = 1446.129
$ws.Range("I136").Value = 1200.88
$ws.Range("K136").Value = 3602.64
$ws.Range("M136").Value = -1052.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 30000
$ws.Range("K33").Value = 30000
$ws.Range("M33").Value = -29748

$ws.Range("H36").Value = 1999.5
$ws.Range("I36").Value = 1999.5
$ws.Range("K36").Value = 1999.5
$ws.Range("M36").Value = -1514.5

$ws.Range("H102").Value = 6922.5386
$ws.Range("I102").Value = 9217
$ws.Range("J102").Value = 3251.4
$ws.Range("K102").Value = 9217
$ws.Range("L102").Value = 3251.4
$ws.Range("M102").Value = -7595
$ws.Range("N102").Value = -6495.4

$ws.Range("H123").Value = 42497.5
$ws.Range("J123").Value = 42497.5
$ws.Range("L123").Value = 42497.5
$ws.Range("N123").Value = -47397.5

$ws.Range("H134").Value = 79995.164
$ws.Range("J134").Value = 79995.164
$ws.Range("L134").Value = 239985.492
$ws.Range("N134").Value = -245055.492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 5000
$ws.Range("K45").Value = 5000
$ws.Range("M45").Value = -4593

$ws.Range("H46").Value = 2448.8333
$ws.Range("J46").Value = 3036.8462
$ws.Range("L46").Value = 3036.8462
$ws.Range("N46").Value = -3412.8462

$ws.Range("H48").Value = 3000
$ws.Range("I48").Value = 3000
$ws.Range("K48").Value = 3000
$ws.Range("M48").Value = -2339

$ws.Range("H105").Value = 41610.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 41610.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 41610.5
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -48598.5

$ws.Range("H122").Value = 4792.394
$ws.Range("I122").Value = 4754.391
$ws.Range("J122").Value = 4879.8
$ws.Range("K122").Value = 14263.173
$ws.Range("L122").Value = 14639.4
$ws.Range("M122").Value = -11813.173
$ws.Range("N122").Value = -19539.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 20175.438
$ws.Range("I107").Value = 1712.2307
$ws.Range("J107").Value = 100182.664
$ws.Range("K107").Value = 5136.6921
$ws.Range("L107").Value = 300547.992
$ws.Range("M107").Value = -3216.6921
$ws.Range("N107").Value = -304387.992

$ws.Range("H122").Value = 3968.5217
$ws.Range("I122").Value = 1989.7241
$ws.Range("J122").Value = 7344.1177
$ws.Range("K122").Value = 5969.1723
$ws.Range("L122").Value = 22032.3531
$ws.Range("M122").Value = -3519.1723
$ws.Range("N122").Value = -26932.3531
